# Add a new "Hình ảnh" (Image) column to the template, with placeholder
# value "[image]" in the sample data row, and grow the sample row's height
# to accommodate an image.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 and placeholder cell F2
$ws.Range("F1").Value = "Hình ảnh"
$ws.Range("F2").Value = "[image]"

# Copy style from existing header/data cells so the new column matches
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122) # xlPasteFormats

# Make the sample row taller to make room for a sample image
$ws.Rows("2").RowHeight = 92.5

# Move the active selection, matching the template's updated default view
$ws.Range("C13").Select()
